# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, identical on both affected sheets.
$updates = @{
    2  = 1077
    5  = 4633
    8  = 1375
    9  = 911
    11 = 1073
    13 = 583
    15 = 12
    16 = 267
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
